$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptos list values (Price + Volume(1h)) to match the latest scrape.
# Some Price cells are plain numeric-looking strings (e.g. "213.58", "2.20") that
# the worksheet stores as TEXT (not Number). Assigning such a string straight to
# .Value would make Excel auto-convert the cell to a Number (and drop things like
# trailing zeros, e.g. "2.20" -> 2.2), so for those cells we force the cell format
# to Text first to keep them as plain text, matching the original data.

# Row 2
$ws.Cells.Item(2, 4).Value = "26.476.65"
$ws.Cells.Item(2, 5).Value = "  -0.95%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.624.65"
$ws.Cells.Item(3, 5).Value = "  +0.00%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.16%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "213.58"
$ws.Cells.Item(5, 5).Value = "  -0.51%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  -0.71%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.19%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.05%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.22%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "19.21"
$ws.Cells.Item(10, 5).Value = "  -0.68%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0856"
$ws.Cells.Item(11, 5).Value = "  -0.11%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "1.853.30"
$ws.Cells.Item(12, 5).Value = "  +0.02%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "1.637.44"
$ws.Cells.Item(13, 5).Value = "  +0.78%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "4.05"
$ws.Cells.Item(14, 5).Value = "  -0.05%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.511"
$ws.Cells.Item(15, 5).Value = "  -0.36%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "64.08"
$ws.Cells.Item(16, 5).Value = "  -1.72%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "234.67"
$ws.Cells.Item(17, 5).Value = "  +0.77%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "26.493.53"
$ws.Cells.Item(18, 5).Value = "  -0.94%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.76"
$ws.Cells.Item(19, 5).Value = "  -0.09%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.25%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.20%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.32"
$ws.Cells.Item(22, 5).Value = "  -1.74%  "

# Row 23
$ws.Cells.Item(23, 2).Value = "Avalanche"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "9.16"
$ws.Cells.Item(23, 5).Value = "  -0.01%  "

# Row 24
$ws.Cells.Item(24, 2).Value = "Toncoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.20"
$ws.Cells.Item(24, 5).Value = "  -1.01%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "146.47"
$ws.Cells.Item(25, 5).Value = "  +0.39%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.14%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "7.07"
$ws.Cells.Item(27, 5).Value = "  +0.34%  "

# Row 28
$ws.Cells.Item(28, 5).Value = "  -0.33%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "15.64"
$ws.Cells.Item(29, 5).Value = "  +0.06%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0496"
$ws.Cells.Item(30, 5).Value = "  -0.64%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.42%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "1.525.30"
$ws.Cells.Item(32, 5).Value = "  +5.08%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "3.27"
$ws.Cells.Item(33, 5).Value = "  +0.43%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.99"
$ws.Cells.Item(34, 5).Value = "  -0.62%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +2.54%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +0.02%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.570"
$ws.Cells.Item(37, 5).Value = "  +0.24%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  -0.78%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.836"
$ws.Cells.Item(39, 5).Value = "  -0.57%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -2.00%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +0.22%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +0.25%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "1.763.93"
$ws.Cells.Item(43, 5).Value = "  -0.01%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "62.92"
$ws.Cells.Item(44, 5).Value = "  +0.96%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.761"
$ws.Cells.Item(45, 5).Value = "  -0.56%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -3.91%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "89.82"
$ws.Cells.Item(47, 5).Value = "  +1.48%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +0.09%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -0.65%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0965"
$ws.Cells.Item(50, 5).Value = "  -0.01%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "7.56"
$ws.Cells.Item(51, 5).Value = "  +0.72%  "

